# Apply a row-level permutation to the data rows (2-17) of the active sheet.
# Each target row ends up containing the full original content (all columns)
# of a specific source row, per the mapping below (target row -> source row).
#
# Columns Y and AA hold date-like text (e.g. "2023-09-03"). Excel's COM
# value-setter auto-parses such strings into date serials, so we temporarily
# force those columns to text format while writing, then restore their
# (default) style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17
$dateTextCols = @("Y", "AA")

foreach ($col in $dateTextCols) {
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).NumberFormat = "@"
}

# Snapshot the full content of every data row BEFORE any writes, so the
# permutation (which includes cycles) can be applied safely.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Rows($r).Value2
}

# Mapping: target row number -> source row number (content to place there)
$mapping = @{
    2  = 9
    3  = 10
    4  = 4
    5  = 7
    6  = 17
    7  = 3
    8  = 16
    9  = 14
    10 = 6
    11 = 2
    12 = 15
    13 = 12
    14 = 11
    15 = 8
    16 = 13
    17 = 5
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $ws.Rows($targetRow).Value2 = $snapshot[$sourceRow]
}

foreach ($col in $dateTextCols) {
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).Style = "Normal"
}
